$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column (27-dec) before the
#     existing "01-oct." column, shifting EX:GB -> EY:GC. ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("EX1").EntireColumn.Insert()
$ws1.Range("EX1").Value = "27-dec"
$ws1.Range("EX2:EX25").Value = "-"

# --- Sheet "Gaz": append a new data row for 2025-12-25. ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A182").Value = "'2025-12-25"
$ws2.Range("B182").Value = 27.5
